$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("table1")
$src = $ws.Range("A6")
Write-Host "src style:" $src.Style
$dst = $ws.Range("A7")
$dst.NumberFormat = $src.NumberFormat
Write-Host "trying style set"
try {
    $dst.Style = $src.Style
    Write-Host "Style set OK"
} catch {
    Write-Host "Style set FAILED:" $_.Exception.Message
}
